$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.732.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.34%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.647.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.46%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'212.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.19%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +3.84%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'23.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.62%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -1.22%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -0.34%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +1.45%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.883.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.28%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.645.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.51%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  -0.64%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.562"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'64.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.24%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'27.667.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.10%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'230.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.59%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "'ShibaInu"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0724"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.74%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "'Chainlink"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'7.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.74%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.08%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.39%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +7.84%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.92%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'149.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.62%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -2.71%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +1.18%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.21%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'15.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.62%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.35%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -2.51%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -0.14%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +1.97%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.438.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.30%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +2.32%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -1.82%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.572"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.12%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.80%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.71%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +14.20%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -1.05%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -0.05%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'5.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.33%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'mCoin"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.09%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'MXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.17%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'65.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.52%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.790.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.32%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.82%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'86.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.01%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0₆0106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.53%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0989"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.32%  "
$ws.Range("E51").Style = "Normal"
